$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New working set of image/word/category cues (row index -> image, word, category).
# Row order matters: writing top-to-bottom lets the shared-string table fill in
# the same sequence the source data was authored in.
$rows = @(
    @('dog/dog007.jpg',    'küssen',   'dog'),
    @('flower/flower022.jpg', 'wecken',   'flower'),
    @('flower/flower000.jpg', 'grenzen',  'flower'),
    @('flower/flower009.jpg', 'legen',    'flower'),
    @('dog/dog001.jpg',    'deuten',   'dog'),
    @('dog/dog030.jpg',    'heißen',   'dog'),
    @('dog/dog020.jpg',    'münzen',   'dog'),
    @('dog/dog025.jpg',    'planen',   'dog'),
    @('flower/flower031.jpg', 'lächeln',  'flower'),
    @('dog/dog022.jpg',    'dienen',   'dog'),
    @('flower/flower025.jpg', 'passen',   'flower'),
    @('dog/dog023.jpg',    'lügen',    'dog'),
    @('dog/dog015.jpg',    'kriegen',  'dog'),
    @('flower/flower007.jpg', 'kranken',  'flower'),
    @('flower/flower026.jpg', 'meinen',   'flower'),
    @('dog/dog016.jpg',    'streifen', 'dog'),
    @('flower/flower029.jpg', 'frischen', 'flower'),
    @('flower/flower021.jpg', 'rufen',    'flower'),
    @('flower/flower002.jpg', 'lernen',   'flower'),
    @('flower/flower014.jpg', 'sparen',   'flower'),
    @('flower/flower003.jpg', 'proben',   'flower'),
    @('dog/dog009.jpg',    'trotzen',  'dog'),
    @('dog/dog026.jpg',    'atmen',    'dog'),
    @('dog/dog027.jpg',    'ändern',   'dog'),
    @('dog/dog013.jpg',    'spüren',   'dog'),
    @('dog/dog014.jpg',    'narren',   'dog'),
    @('flower/flower010.jpg', 'stoßen',   'flower'),
    @('flower/flower019.jpg', 'danken',   'flower'),
    @('flower/flower027.jpg', 'parken',   'flower'),
    @('dog/dog028.jpg',    'hacken',   'dog'),
    @('dog/dog000.jpg',    'lassen',   'dog'),
    @('flower/flower020.jpg', 'angeln',   'flower')
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $r++
}
